$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.4483165197934
$ws.Range("C2").Value = 9.154983712916174
$ws.Range("D2").Value = 3.615603583881652
$ws.Range("F2").Value = 17.04688280068747
$ws.Range("G2").Value = 16.43456771822514
$ws.Range("H2").Value = 11.08967466736438
$ws.Range("N2").Value = 15.50902330551442
$ws.Range("O2").Value = 15.22358615560785
$ws.Range("B3").Value = 10.82451622253506
$ws.Range("C3").Value = 8.833332285859225
$ws.Range("D3").Value = 3.508930007559171
$ws.Range("F3").Value = 17.04300878553376
$ws.Range("G3").Value = 16.39678749686747
$ws.Range("H3").Value = 11.13180457383454
$ws.Range("N3").Value = 15.51592581130319
$ws.Range("O3").Value = 15.28006549916674
$ws.Range("B4").Value = 10.42239480811681
$ws.Range("C4").Value = 8.628847687521185
$ws.Range("D4").Value = 3.441098456653105
$ws.Range("F4").Value = 17.04681965753492
$ws.Range("G4").Value = 16.38248584184311
$ws.Range("H4").Value = 11.15981399388763
$ws.Range("N4").Value = 15.52186180946102
$ws.Range("O4").Value = 15.31921114959527
$ws.Range("B5").Value = 10.25386068586346
$ws.Range("C5").Value = 8.543855826896881
$ws.Range("D5").Value = 3.412894884023953
$ws.Range("F5").Value = 17.04992805078433
$ws.Range("G5").Value = 16.37889595074702
$ws.Range("H5").Value = 11.1717660535306
$ws.Range("N5").Value = 15.52470902946096
$ws.Range("O5").Value = 15.33628228711098
$ws.Range("B6").Value = 10.22559846443014
$ws.Range("C6").Value = 8.529645539883514
$ws.Range("D6").Value = 3.408178557571687
$ws.Range("F6").Value = 17.05053808358255
$ws.Range("G6").Value = 16.37843498880328
$ws.Range("H6").Value = 11.17378316273779
$ws.Range("N6").Value = 15.52520771090922
$ws.Range("O6").Value = 15.33918438915771
$ws.Range("B7").Value = 10.42014058953104
$ws.Range("C7").Value = 8.627708057124375
$ws.Range("D7").Value = 3.440720332726455
$ws.Range("F7").Value = 17.04685528292561
$ws.Range("G7").Value = 16.38242836644552
$ws.Range("H7").Value = 11.15997300589254
$ws.Range("N7").Value = 15.52189847239062
$ws.Range("O7").Value = 15.31943685212928
$ws.Range("B8").Value = 11.23727001136684
$ws.Range("C8").Value = 9.045582678014313
$ws.Range("D8").Value = 3.579321238480738
$ws.Range("F8").Value = 17.04426226583156
$ws.Range("G8").Value = 16.41969684952241
$ws.Range("H8").Value = 11.10375629489078
$ws.Range("N8").Value = 15.51105171162833
$ws.Range("O8").Value = 15.24213054597412
$ws.Range("B9").Value = 12.68327450934834
$ws.Range("C9").Value = 9.805916249799466
$ws.Range("D9").Value = 3.831611477533428
$ws.Range("F9").Value = 17.08826323817607
$ws.Range("G9").Value = 16.56310556871167
$ws.Range("H9").Value = 11.01052936186318
$ws.Range("N9").Value = 15.50319376337328
$ws.Range("O9").Value = 15.12616218188715
$ws.Range("B10").Value = 13.64535287856902
$ws.Range("C10").Value = 10.32425028905343
$ws.Range("D10").Value = 4.003904411565509
$ws.Range("F10").Value = 17.15037999818984
$ws.Range("G10").Value = 16.71070760611783
$ws.Range("H10").Value = 10.95243750530239
$ws.Range("N10").Value = 15.50551625554769
$ws.Range("O10").Value = 15.06291963980145
$ws.Range("B11").Value = 14.063719302182
$ws.Range("C11").Value = 10.55058098332583
$ws.Range("D11").Value = 4.079240451852024
$ws.Range("F11").Value = 17.18504502969394
$ws.Range("G11").Value = 16.78680729261561
$ws.Range("H11").Value = 10.92827506371955
$ws.Range("N11").Value = 15.50831220423471
$ws.Range("O11").Value = 15.03896342358866
$ws.Range("B12").Value = 14.22800221528852
$ws.Range("C12").Value = 10.63487538862329
$ws.Range("D12").Value = 4.107316081304415
$ws.Range("F12").Value = 17.19908545392696
$ws.Range("G12").Value = 16.8168877106194
$ws.Range("H12").Value = 10.91945146434235
$ws.Range("N12").Value = 15.50961937079761
$ws.Range("O12").Value = 15.03058731696572
$ws.Range("B13").Value = 14.1927950636679
$ws.Range("C13").Value = 10.6167846326398
$ws.Range("D13").Value = 4.101289841150016
$ws.Range("F13").Value = 17.19602111056205
$ws.Range("G13").Value = 16.81035360642258
$ws.Range("H13").Value = 10.92133726647455
$ws.Range("N13").Value = 15.50932682966961
$ws.Range("O13").Value = 15.03236027265626
$ws.Range("B14").Value = 14.07731342454765
$ws.Range("C14").Value = 10.55754450658238
$ws.Range("D14").Value = 4.081559400535938
$ws.Range("F14").Value = 17.18618188625416
$ws.Range("G14").Value = 16.78925686550901
$ws.Range("H14").Value = 10.92754259977066
$ws.Range("N14").Value = 15.50841477743901
$ws.Range("O14").Value = 15.03826035672883
$ws.Range("B15").Value = 14.00655676811129
$ws.Range("C15").Value = 10.52107290736473
$ws.Range("D15").Value = 4.069414585456664
$ws.Range("F15").Value = 17.18027378412612
$ws.Range("G15").Value = 16.7764982067977
$ws.Range("H15").Value = 10.93138604485345
$ws.Range("N15").Value = 15.50788841577512
$ws.Range("O15").Value = 15.04196501317512
$ws.Range("B16").Value = 13.61776360062831
$ws.Range("C16").Value = 10.30926408736229
$ws.Range("D16").Value = 3.998918493405121
$ws.Range("F16").Value = 17.14824285870548
$ws.Range("G16").Value = 16.70591248911777
$ws.Range("H16").Value = 10.95406218811643
$ws.Range("N16").Value = 15.50536838968022
$ws.Range("O16").Value = 15.0645823907657
$ws.Range("B17").Value = 13.37345949586183
$ws.Range("C17").Value = 10.17686394410599
$ws.Range("D17").Value = 3.954881200174471
$ws.Range("F17").Value = 17.13022899594972
$ws.Range("G17").Value = 16.66488766185379
$ws.Range("H17").Value = 10.96855355951599
$ws.Range("N17").Value = 15.50426674134384
$ws.Range("O17").Value = 15.07969274014799
$ws.Range("B18").Value = 13.23082870127223
$ws.Range("C18").Value = 10.09982288691806
$ws.Range("D18").Value = 3.929266687212492
$ws.Range("F18").Value = 17.12047160644269
$ws.Range("G18").Value = 16.64213584384553
$ws.Range("H18").Value = 10.97710165721162
$ws.Range("N18").Value = 15.50379693430375
$ws.Range("O18").Value = 15.08883665902461
$ws.Range("B19").Value = 13.18217491437188
$ws.Range("C19").Value = 10.07358727616294
$ws.Range("D19").Value = 3.92054552312915
$ws.Range("F19").Value = 17.11727182430051
$ws.Range("G19").Value = 16.63457821457355
$ws.Range("H19").Value = 10.9800324666302
$ws.Range("N19").Value = 15.50366605540701
$ws.Range("O19").Value = 15.09201029897292
$ws.Range("B20").Value = 13.39968518821565
$ws.Range("C20").Value = 10.19105046428342
$ws.Range("D20").Value = 3.959598701505541
$ws.Range("F20").Value = 17.13208417312732
$ws.Range("G20").Value = 16.6691675805159
$ws.Range("H20").Value = 10.96698887424038
$ws.Range("N20").Value = 15.50436706997098
$ws.Range("O20").Value = 15.07803732377067
$ws.Range("B21").Value = 14.11133943980725
$ws.Range("C21").Value = 10.57498346884958
$ws.Range("D21").Value = 4.087367102587248
$ws.Range("F21").Value = 17.18904718285208
$ws.Range("G21").Value = 16.79541942275914
$ws.Range("H21").Value = 10.92571108499614
$ws.Range("N21").Value = 15.5086759418803
$ws.Range("O21").Value = 15.03650845465905
$ws.Range("B22").Value = 14.58223370353113
$ws.Range("C22").Value = 10.81765573754221
$ws.Range("D22").Value = 4.168227778038489
$ws.Range("F22").Value = 17.23159577580296
$ws.Range("G22").Value = 16.88528082575252
$ws.Range("H22").Value = 10.90063526315705
$ws.Range("N22").Value = 15.51293905834316
$ws.Range("O22").Value = 15.01342277996925
$ws.Range("B23").Value = 14.33299524546446
$ws.Range("C23").Value = 10.68890716136954
$ws.Range("D23").Value = 4.125317373324981
$ws.Range("F23").Value = 17.20840303093636
$ws.Range("G23").Value = 16.8366567517059
$ws.Range("H23").Value = 10.91384449024564
$ws.Range("N23").Value = 15.51053194049627
$ws.Range("O23").Value = 15.02537186052055
$ws.Range("B24").Value = 13.38783533038008
$ws.Range("C24").Value = 10.18463960778451
$ws.Range("D24").Value = 3.957466842721239
$ws.Range("F24").Value = 17.1312435810609
$ws.Range("G24").Value = 16.66723002995779
$ws.Range("H24").Value = 10.96769559313823
$ws.Range("N24").Value = 15.50432120183923
$ws.Range("O24").Value = 15.07878431495727
$ws.Range("B25").Value = 12.30941591217682
$ws.Range("C25").Value = 9.607027169801835
$ws.Range("D25").Value = 3.765577688424414
$ws.Range("F25").Value = 17.07111246903078
$ws.Range("G25").Value = 16.51683005777631
$ws.Range("H25").Value = 11.03392516510445
$ws.Range("N25").Value = 15.503891715315
$ws.Range("O25").Value = 15.15369473665177
